$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.774.77'
$ws.Range('E2').Value = '  +3.25%  '
$ws.Range('D3').Value = '3.135.16'
$ws.Range('E3').Value = '  +2.06%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('E5').Value = '  +1.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.81'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +3.64%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.129.38'
$ws.Range('E8').Value = '  +2.27%  '
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('E10').Value = '  +16.15%  '
$ws.Range('E11').Value = '  +2.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.469'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('E13').Value = '  +5.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.64'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +4.10%  '
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').Value = '3.653.30'
$ws.Range('E16').Value = '  +2.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.18'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').Value = '63.689.34'
$ws.Range('E18').Value = '  +3.23%  '
$ws.Range('D19').Value = '3.129.74'
$ws.Range('E19').Value = '  +1.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '465.66'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +3.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.41'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +3.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.735'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.30'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -3.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.50'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.96'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +10.88%  '
$ws.Range('E28').Value = '  +2.38%  '
$ws.Range('E29').Value = '  -0.92%  '
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  +2.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.20'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +2.05%  '
$ws.Range('E33').Value = '  -3.44%  '
$ws.Range('D34').Value = '0.0₃0869'
$ws.Range('E34').Value = '  +9.31%  '
$ws.Range('E35').Value = '  +9.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.06'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +1.85%  '
$ws.Range('E37').Value = '  +13.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.12'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +1.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '51.05'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +1.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '447.48'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +4.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.78'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0373'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +0.45%  '
$ws.Range('D43').Value = '2.904.21'
$ws.Range('E43').Value = '  +3.76%  '
$ws.Range('E44').Value = '  +4.08%  '
$ws.Range('E45').Value = '  +2.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.20'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +5.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.18'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +3.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.95'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +0.75%  '
$ws.Range('E51').Value = '  +3.87%  '

Write-Host "Applied cryptos update"
